# edit.ps1
# Applies the weekly crime-data refresh described in the commit:
# "New crime data collected"
#
# 1) Update the "Volume 30 Number NN" label (17 -> 18)
# 2) Update the reporting week range (4/24/2023-4/30/2023 -> 5/1/2023-5/7/2023)
# 3) Update every weekly/28-day/YTD/2-year crime-count and %-change figure
#    in the main CompStat table (rows 14-30) to the newly collected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------

# A8 holds rich text: "Volume 30   Number  17" -> replace the trailing "17" with "18"
$a8 = $ws.Range("A8")
$a8text = $a8.Value2
$idx = $a8text.LastIndexOf("17")
$a8.Characters($idx + 1, 2).Text = "18"

# C9 holds rich text: "Report Covering the Week  4/24/2023  Through  4/30/2023"
# -> update both date runs to the new reporting week
$c9 = $ws.Range("C9")
$c9text = $c9.Value2
$idx1 = $c9text.IndexOf("4/24/2023")
$c9.Characters($idx1 + 1, 9).Text = "5/1/2023"
$c9text = $c9.Value2
$idx2 = $c9text.IndexOf("4/30/2023")
$c9.Characters($idx2 + 1, 9).Text = "5/7/2023"

# --- CompStat data table updates (rows 14-30) -----------------------------

$ws.Range("D14").Value = 2
$ws.Range("G14").Value = 5
$ws.Range("J14").Value = 15
$ws.Range("K14").Value = -80
$ws.Range("N14").Value = -91.428571428571
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 14
$ws.Range("H15").Value = -42.857142857142
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = -43.333333333333
$ws.Range("L15").Value = -34.615384615384
$ws.Range("M15").Value = 6.25
$ws.Range("N15").Value = -48.484848484848
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = 45
$ws.Range("E16").Value = -15.555555555555
$ws.Range("F16").Value = 129
$ws.Range("G16").Value = 160
$ws.Range("H16").Value = -19.375
$ws.Range("I16").Value = 612
$ws.Range("J16").Value = 694
$ws.Range("K16").Value = -11.815561959654
$ws.Range("L16").Value = 50.738916256157
$ws.Range("M16").Value = 53
$ws.Range("N16").Value = -82.255726297477
$ws.Range("C17").Value = 33
$ws.Range("D17").Value = 47
$ws.Range("E17").Value = -29.787234042553
$ws.Range("F17").Value = 148
$ws.Range("G17").Value = 166
$ws.Range("H17").Value = -10.843373493975
$ws.Range("I17").Value = 719
$ws.Range("J17").Value = 633
$ws.Range("K17").Value = 13.586097946287
$ws.Range("L17").Value = 29.084380610412
$ws.Range("M17").Value = 66.050808314087
$ws.Range("N17").Value = -36.032028469750
$ws.Range("C18").Value = 45
$ws.Range("D18").Value = 62
$ws.Range("E18").Value = -27.419354838709
$ws.Range("F18").Value = 178
$ws.Range("G18").Value = 225
$ws.Range("H18").Value = -20.888888888888
$ws.Range("I18").Value = 735
$ws.Range("J18").Value = 987
$ws.Range("K18").Value = -25.531914893617
$ws.Range("L18").Value = 13.425925925925
$ws.Range("M18").Value = 10.526315789473
$ws.Range("N18").Value = -81.865284974093
$ws.Range("C19").Value = 228
$ws.Range("D19").Value = 232
$ws.Range("E19").Value = -1.724137931034
$ws.Range("F19").Value = 913
$ws.Range("G19").Value = 848
$ws.Range("H19").Value = 7.665094339622
$ws.Range("I19").Value = 3884
$ws.Range("J19").Value = 3592
$ws.Range("K19").Value = 8.129175946547
$ws.Range("L19").Value = 88.543689320388
$ws.Range("M19").Value = 11.769784172661
$ws.Range("N19").Value = -63.926813411349
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -15.384615384615
$ws.Range("F20").Value = 55
$ws.Range("G20").Value = 48
$ws.Range("H20").Value = 14.583333333333
$ws.Range("I20").Value = 188
$ws.Range("J20").Value = 177
$ws.Range("K20").Value = 6.214689265536
$ws.Range("L20").Value = 34.285714285714
$ws.Range("M20").Value = 62.068965517241
$ws.Range("N20").Value = -91.304347826087
$ws.Range("C21").Value = 356
$ws.Range("D21").Value = 404
$ws.Range("E21").Value = -11.881188118811
$ws.Range("F21").Value = 1431
$ws.Range("G21").Value = 1466
$ws.Range("H21").Value = -2.387448840381
$ws.Range("I21").Value = 6175
$ws.Range("J21").Value = 6158
$ws.Range("K21").Value = 0.276063657031
$ws.Range("L21").Value = 59.560723514211
$ws.Range("M21").Value = 20.417316692667
$ws.Range("N21").Value = -71.485962319911
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 20
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 45
$ws.Range("G22").Value = 61
$ws.Range("H22").Value = -26.229508196721
$ws.Range("I22").Value = 216
$ws.Range("J22").Value = 247
$ws.Range("K22").Value = -12.550607287449
$ws.Range("L22").Value = 45.945945945945
$ws.Range("M22").Value = 14.893617021276
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("G23").Value = 38
$ws.Range("H23").Value = -21.052631578947
$ws.Range("I23").Value = 131
$ws.Range("J23").Value = 157
$ws.Range("K23").Value = -16.560509554140
$ws.Range("L23").Value = -25.988700564971
$ws.Range("M23").Value = 1.550387596899
$ws.Range("C24").Value = 380
$ws.Range("D24").Value = 414
$ws.Range("E24").Value = -8.212560386473
$ws.Range("F24").Value = 1599
$ws.Range("G24").Value = 1686
$ws.Range("H24").Value = -5.160142348754
$ws.Range("I24").Value = 6563
$ws.Range("J24").Value = 6902
$ws.Range("K24").Value = -4.911619820341
$ws.Range("L24").Value = 56.934481109517
$ws.Range("M24").Value = 18.916470375067
$ws.Range("C25").Value = 92
$ws.Range("D25").Value = 101
$ws.Range("E25").Value = -8.910891089108
$ws.Range("G25").Value = 375
$ws.Range("H25").Value = 2.666666666666
$ws.Range("I25").Value = 1609
$ws.Range("J25").Value = 1502
$ws.Range("K25").Value = 7.123834886817
$ws.Range("L25").Value = 49.674418604651
$ws.Range("M25").Value = 28.207171314741
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 73
$ws.Range("J26").Value = 92
$ws.Range("K26").Value = -20.652173913043
$ws.Range("L26").Value = -9.876543209876
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 21
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 78
$ws.Range("G27").Value = 84
$ws.Range("H27").Value = -7.142857142857
$ws.Range("I27").Value = 305
$ws.Range("J27").Value = 303
$ws.Range("K27").Value = 0.660066006600
$ws.Range("L27").Value = 39.908256880733
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -35
$ws.Range("L28").Value = 8.333333333333
$ws.Range("M28").Value = 8.333333333333
$ws.Range("N28").Value = -77.192982456140
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 11
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = -35.294117647058
$ws.Range("L29").Value = -8.333333333333
$ws.Range("M29").Value = 37.5
$ws.Range("N29").Value = -78
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = -18.181818181818
$ws.Range("I30").Value = 34
$ws.Range("J30").Value = 71
$ws.Range("K30").Value = -52.112676056338
$ws.Range("L30").Value = -30.612244897959
